$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 2P": fill in real Segundo Parcial grade stats ---
$ws2P = $wb.Worksheets.Item("Estadisticos 2P")
$ws2P.Range("D2").Value = 0
$ws2P.Range("E2").Value = 3
$ws2P.Range("F2").Value = 28
$ws2P.Range("G2").Value = 90.31999999999999
$ws2P.Range("H2").Value = 7.7

# --- Sheet "Estadisticos Final": recalculated final stats ---
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")
$wsFinal.Range("E2").Value = 3
$wsFinal.Range("F2").Value = 28
$wsFinal.Range("G2").Value = 90.31999999999999
$wsFinal.Range("H2").Value = 7.9

# --- Sheet "Rescatables": only students still needing retake remain ---
$wsResc = $wb.Worksheets.Item("Rescatables")
# Delete from the bottom up so row numbers of earlier rows stay valid.
$wsResc.Rows.Item(7).Delete()
$wsResc.Rows.Item(6).Delete()
$wsResc.Rows.Item(4).Delete()
$wsResc.Rows.Item(2).Delete()
# Update remaining "Reprobadas" count for the first student (4 -> 3).
$wsResc.Range("G2").Value = 3
